$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ten_lists")

# Update the list identifier label ("S000" -> "S017 OA")
$ws.Range("F1").Value = "S017 OA"

# The walk/no-walk condition labels got reshuffled between a handful of cells
$ws.Range("I3").Value = "no walk/diff"
$ws.Range("C10").Value = "no walk/same"
$ws.Range("C17").Value = "no walk/diff"
$ws.Range("I17").Value = "walk/same"
$ws.Range("C24").Value = "walk/diff"
$ws.Range("C31").Value = "walk/diff"

# Leave the selection where the author last left it
$ws.Range("I24").Select() | Out-Null
